$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase de Grupos")
$ws.Unprotect("CC01")
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("K16").Formula = $ws.Range("K16").Formula
$ws.Range("L16").Formula = $ws.Range("L16").Formula
$ws.Protect("CC01")
